$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2021_bottom4")
$ws.Cells.Item(2, 2).Value = 0.5998
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(2, 4).Value = 0.3002
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(3, 2).Value = 0.85
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05

$ws = $wb.Worksheets.Item("2020_top8")
$ws.Cells.Item(2, 2).Value = 0.2536
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(2, 4).Value = 0.161
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(2, 6).Value = 0.0939
$ws.Cells.Item(2, 7).Value = 0.05
$ws.Cells.Item(2, 8).Value = 0.1631
$ws.Cells.Item(2, 9).Value = 0.1784
$ws.Cells.Item(3, 2).Value = 0.05
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05
$ws.Cells.Item(3, 6).Value = 0.05
$ws.Cells.Item(3, 7).Value = 0.65
$ws.Cells.Item(3, 8).Value = 0.05
$ws.Cells.Item(3, 9).Value = 0.05

$ws = $wb.Worksheets.Item("2022_top8")
$ws.Cells.Item(2, 2).Value = 0.3185
$ws.Cells.Item(2, 3).Value = 0.0659
$ws.Cells.Item(2, 4).Value = 0.1297
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(2, 6).Value = 0.05
$ws.Cells.Item(2, 7).Value = 0.05
$ws.Cells.Item(2, 8).Value = 0.2859
$ws.Cells.Item(2, 9).Value = 0.05
$ws.Cells.Item(3, 2).Value = 0.4022
$ws.Cells.Item(3, 3).Value = 0.2978
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05
$ws.Cells.Item(3, 6).Value = 0.05
$ws.Cells.Item(3, 7).Value = 0.05
$ws.Cells.Item(3, 8).Value = 0.05
$ws.Cells.Item(3, 9).Value = 0.05

$ws = $wb.Worksheets.Item("2021_top8")
$ws.Cells.Item(2, 2).Value = 0.05
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(2, 4).Value = 0.05
$ws.Cells.Item(2, 5).Value = 0.1554
$ws.Cells.Item(2, 6).Value = 0.3822
$ws.Cells.Item(2, 7).Value = 0.064
$ws.Cells.Item(2, 8).Value = 0.0544
$ws.Cells.Item(2, 9).Value = 0.194
$ws.Cells.Item(3, 2).Value = 0.3097
$ws.Cells.Item(3, 3).Value = 0.2894
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05
$ws.Cells.Item(3, 6).Value = 0.151
$ws.Cells.Item(3, 7).Value = 0.05
$ws.Cells.Item(3, 8).Value = 0.05
$ws.Cells.Item(3, 9).Value = 0.05

$ws = $wb.Worksheets.Item("maxSRWeights")
$ws.Cells.Item(2, 2).Value = 0.05
$ws.Cells.Item(2, 3).Value = 0.0500000160251316
$ws.Cells.Item(2, 4).Value = 0.05
$ws.Cells.Item(3, 2).Value = 0.05
$ws.Cells.Item(3, 3).Value = 0.0499999999999998
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(4, 2).Value = 0.0500000694278713
$ws.Cells.Item(4, 3).Value = 0.05
$ws.Cells.Item(4, 4).Value = 0.05
$ws.Cells.Item(5, 2).Value = 0.0500000000000001
$ws.Cells.Item(5, 3).Value = 0.05
$ws.Cells.Item(5, 4).Value = 0.05
$ws.Cells.Item(6, 2).Value = 0.0500000000000003
$ws.Cells.Item(6, 3).Value = 0.05
$ws.Cells.Item(6, 4).Value = 0.05
$ws.Cells.Item(7, 2).Value = 0.05
$ws.Cells.Item(7, 3).Value = 0.05
$ws.Cells.Item(7, 4).Value = 0.0629489401567271
$ws.Cells.Item(8, 2).Value = 0.05
$ws.Cells.Item(8, 3).Value = 0.0500000000000001
$ws.Cells.Item(8, 4).Value = 0.05
$ws.Cells.Item(9, 2).Value = 0.0499999999999999
$ws.Cells.Item(9, 3).Value = 0.0500000000000001
$ws.Cells.Item(9, 4).Value = 0.0500000000000001
$ws.Cells.Item(10, 2).Value = 0.0499999999999999
$ws.Cells.Item(10, 3).Value = 0.0499999999999999
$ws.Cells.Item(10, 4).Value = 0.05
$ws.Cells.Item(11, 2).Value = 0.0999999305721293
$ws.Cells.Item(11, 3).Value = 0.09999998397486889
$ws.Cells.Item(11, 4).Value = 0.05
$ws.Cells.Item(12, 2).Value = 0.05
$ws.Cells.Item(12, 3).Value = 0.0500000000000001
$ws.Cells.Item(12, 4).Value = 0.0499999999999999
$ws.Cells.Item(13, 2).Value = 0.0500000000000001
$ws.Cells.Item(13, 3).Value = 0.05
$ws.Cells.Item(13, 4).Value = 0.05
$ws.Cells.Item(14, 2).Value = 0.05
$ws.Cells.Item(14, 3).Value = 0.05
$ws.Cells.Item(14, 4).Value = 0.05
$ws.Cells.Item(15, 2).Value = 0.05
$ws.Cells.Item(15, 3).Value = 0.05
$ws.Cells.Item(15, 4).Value = 0.05
$ws.Cells.Item(16, 2).Value = 0.05
$ws.Cells.Item(16, 3).Value = 0.05
$ws.Cells.Item(16, 4).Value = 0.0870510598432723
$ws.Cells.Item(17, 2).Value = 0.05
$ws.Cells.Item(17, 3).Value = 0.05
$ws.Cells.Item(17, 4).Value = 0.0500000000000002
$ws.Cells.Item(18, 2).Value = 0.0500000000000001
$ws.Cells.Item(18, 3).Value = 0.05
$ws.Cells.Item(18, 4).Value = 0.05
$ws.Cells.Item(19, 2).Value = 0.0500000000000001
$ws.Cells.Item(19, 3).Value = 0.05
$ws.Cells.Item(19, 4).Value = 0.0500000000000001
$ws.Cells.Item(20, 2).Value = 0.05
$ws.Cells.Item(20, 3).Value = 0.05
$ws.Cells.Item(20, 4).Value = 0.05

$ws = $wb.Worksheets.Item("MVPWeights")
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(2, 4).Value = 0.1
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(4, 3).Value = 0.05
$ws.Cells.Item(4, 4).Value = 0.05
$ws.Cells.Item(5, 3).Value = 0.05
$ws.Cells.Item(5, 4).Value = 0.05
$ws.Cells.Item(6, 3).Value = 0.05
$ws.Cells.Item(6, 4).Value = 0.05
$ws.Cells.Item(7, 3).Value = 0.0499999999999996
$ws.Cells.Item(7, 4).Value = 0.05
$ws.Cells.Item(8, 3).Value = 0.05
$ws.Cells.Item(8, 4).Value = 0.05
$ws.Cells.Item(9, 3).Value = 0.0500000000000001
$ws.Cells.Item(9, 4).Value = 0.0500000000000001
$ws.Cells.Item(10, 3).Value = 0.0500000000000001
$ws.Cells.Item(10, 4).Value = 0.0499999999999999
$ws.Cells.Item(11, 3).Value = 0.05
$ws.Cells.Item(11, 4).Value = 0.05
$ws.Cells.Item(12, 3).Value = 0.05
$ws.Cells.Item(12, 4).Value = 0.05
$ws.Cells.Item(13, 3).Value = 0.0499999999999999
$ws.Cells.Item(13, 4).Value = 0.0499999999999999
$ws.Cells.Item(14, 3).Value = 0.05
$ws.Cells.Item(14, 4).Value = 0.05
$ws.Cells.Item(15, 3).Value = 0.05
$ws.Cells.Item(15, 4).Value = 0.0500000000000003
$ws.Cells.Item(16, 3).Value = 0.05
$ws.Cells.Item(16, 4).Value = 0.05
$ws.Cells.Item(17, 3).Value = 0.0500000000000001
$ws.Cells.Item(17, 4).Value = 0.0500000000000001
$ws.Cells.Item(18, 3).Value = 0.09999999999999989
$ws.Cells.Item(18, 4).Value = 0.05
$ws.Cells.Item(19, 3).Value = 0.05
$ws.Cells.Item(19, 4).Value = 0.05
$ws.Cells.Item(20, 3).Value = 0.05
$ws.Cells.Item(20, 4).Value = 0.05

$ws = $wb.Worksheets.Item("2020_top4")
$ws.Cells.Item(2, 2).Value = 0.1634
$ws.Cells.Item(2, 3).Value = 0.6529
$ws.Cells.Item(2, 4).Value = 0.1338
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(3, 2).Value = 0.05
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.85

$ws = $wb.Worksheets.Item("2022_top4")
$ws.Cells.Item(2, 2).Value = 0.2707
$ws.Cells.Item(2, 3).Value = 0.1503
$ws.Cells.Item(2, 4).Value = 0.529
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(3, 2).Value = 0.85
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05

$ws = $wb.Worksheets.Item("2021_top4")
$ws.Cells.Item(3, 2).Value = 0.85
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05

$ws = $wb.Worksheets.Item("2022_bottom4")
$ws.Cells.Item(3, 2).Value = 0.05
$ws.Cells.Item(3, 3).Value = 0.3058
$ws.Cells.Item(3, 4).Value = 0.4692
$ws.Cells.Item(3, 5).Value = 0.1751

$ws = $wb.Worksheets.Item("2020_bottom4")
$ws.Cells.Item(2, 2).Value = 0.85
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(2, 4).Value = 0.05
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(3, 2).Value = 0.05
$ws.Cells.Item(3, 3).Value = 0.85
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05

$ws = $wb.Worksheets.Item("2022_bottom8")
$ws.Cells.Item(2, 2).Value = 0.65
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(2, 4).Value = 0.05
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(2, 6).Value = 0.05
$ws.Cells.Item(2, 7).Value = 0.05
$ws.Cells.Item(2, 8).Value = 0.05
$ws.Cells.Item(2, 9).Value = 0.05
$ws.Cells.Item(3, 2).Value = 0.05
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.2146
$ws.Cells.Item(3, 5).Value = 0.05
$ws.Cells.Item(3, 6).Value = 0.05
$ws.Cells.Item(3, 7).Value = 0.3215
$ws.Cells.Item(3, 8).Value = 0.2138
$ws.Cells.Item(3, 9).Value = 0.05

$ws = $wb.Worksheets.Item("2020_bottom8")
$ws.Cells.Item(2, 2).Value = 0.65
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(2, 4).Value = 0.05
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(2, 6).Value = 0.05
$ws.Cells.Item(2, 7).Value = 0.05
$ws.Cells.Item(2, 8).Value = 0.05
$ws.Cells.Item(2, 9).Value = 0.05
$ws.Cells.Item(3, 2).Value = 0.05
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05
$ws.Cells.Item(3, 6).Value = 0.05
$ws.Cells.Item(3, 7).Value = 0.65
$ws.Cells.Item(3, 8).Value = 0.05
$ws.Cells.Item(3, 9).Value = 0.05

$ws = $wb.Worksheets.Item("2021_bottom8")
$ws.Cells.Item(2, 2).Value = 0.3383
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(2, 4).Value = 0.1939
$ws.Cells.Item(2, 5).Value = 0.05
$ws.Cells.Item(2, 6).Value = 0.05
$ws.Cells.Item(2, 7).Value = 0.05
$ws.Cells.Item(2, 8).Value = 0.05
$ws.Cells.Item(2, 9).Value = 0.2178
$ws.Cells.Item(3, 2).Value = 0.2991
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(3, 4).Value = 0.05
$ws.Cells.Item(3, 5).Value = 0.05
$ws.Cells.Item(3, 6).Value = 0.4009
$ws.Cells.Item(3, 7).Value = 0.05
$ws.Cells.Item(3, 8).Value = 0.05
$ws.Cells.Item(3, 9).Value = 0.05
